$d = $word.ActiveDocument

# 1. Remove the "Meta description: ..." paragraph that currently sits right
#    after the title heading (paragraph 2).
$metaPara = $d.Paragraphs(2)
$metaPara.Range.Delete()

# 2. At the very end of the document, the last paragraph currently holds the
#    italic AI image-prompt text. Replace it with two paragraphs:
#      - a new bold paragraph repeating the page title
#      - the same italic paragraph, but now containing the meta-description text
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$r = $lastPara.Range
$r.Collapse(1)

$xml = @"
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Coin-O-Mania for Free: Review of IGT's Treasure Hunt Adventure Slot</w:t></w:r></w:p><w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Join a crew on a treasure hunt in Coin-O-Mania, the 5-reel, 4-row slot game from IGT. Play for free and read our review for more information.</w:t></w:r></w:p>
"@
$r.InsertXML($xml)
